$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 877
$ws.Range("B2").Value = 874
$ws.Range("C2").Value = 874
$ws.Range("D2").Value = 874
$ws.Range("E2").Value = 877
$ws.Range("F2").Value = 1000
$ws.Range("G2").Value = 854
$ws.Range("H2").Value = 877
